$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1673.93
$ws.Range("C3").Value = 2324.07
$ws.Range("C4").Value = 932.41
$ws.Range("C9").Value = 47.62
